$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Spencer added hours worked on Thursday (G5)
$ws.Range("G5").Value = 4

# Update the selection/active cell to match the author's final cursor position
$ws.Range("G6").Select()
